$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    $r.End = $r.End - 1
    $r.Text = $newText
}

Set-CellText $t 1 1 "370×7=2590"
Set-CellText $t 1 2 "634×5=3170"
Set-CellText $t 1 3 "407×6=2442"
Set-CellText $t 1 4 "424×6=2544"
Set-CellText $t 1 5 "178×8=1424"

Set-CellText $t 5 1 "675×2=1350"
Set-CellText $t 5 2 "554×2=1108"
Set-CellText $t 5 3 "374×7=2618"
Set-CellText $t 5 4 "313×6=1878"
Set-CellText $t 5 5 "904×3=2712"

Set-CellText $t 10 1 "687×3=2061"
Set-CellText $t 10 2 "919×8=7352"
Set-CellText $t 10 3 "446×6=2676"
Set-CellText $t 10 4 "225×9=2025"
Set-CellText $t 10 5 "879×8=7032"

Set-CellText $t 15 1 "113×7=791"
Set-CellText $t 15 2 "634×2=1268"
Set-CellText $t 15 3 "499×8=3992"
Set-CellText $t 15 4 "249×4=996"
Set-CellText $t 15 5 "384×9=3456"

Set-CellText $t 20 1 "233×9=2097"
Set-CellText $t 20 2 "714×8=5712"
Set-CellText $t 20 3 "941×2=1882"
Set-CellText $t 20 4 "431×9=3879"
Set-CellText $t 20 5 "707×2=1414"
